$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Remove the standalone "Extend table as required" paragraph that
#    sits just before the "Function Declarations" Heading 1 paragraph.
#    (The paragraph's whole range -- including its paragraph mark -- is
#    deleted so the following Heading 1 paragraph is unaffected.)
# ---------------------------------------------------------------------
$extendTable1 = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "Extend table as required`r") {
        $extendTable1 = $p
    }
}
if ($extendTable1 -ne $null) {
    $extendTable1.Range.Delete()
}

# ---------------------------------------------------------------------
# 2) Clear the "Only include functions that you will develop." run from
#    its (italic) paragraph, leaving the paragraph mark/formatting in
#    place, then add a second, identically formatted, empty paragraph
#    right after it.
# ---------------------------------------------------------------------
$find1 = $d.Content
$find1.Find.ClearFormatting()
$find1.Find.Replacement.ClearFormatting()
$find1.Find.Execute("Only include functions that you will develop.", $true, $false, $false, $false, $false, $true, 1, $false, "^p", 2) | Out-Null

# ---------------------------------------------------------------------
# 3) Remove the italic "Extend table as required. Note that..." paragraph
#    right before the "Flowchart(s)" Heading 1 paragraph (near the end
#    of the document).
# ---------------------------------------------------------------------
$extendTable2 = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Extend table as required. Note that*") {
        $extendTable2 = $p
    }
}
if ($extendTable2 -ne $null) {
    $extendTable2.Range.Delete()
}

# ---------------------------------------------------------------------
# 4) Reword "May be included as separate pdf" -> "Included as separate PDF"
# ---------------------------------------------------------------------
$find2 = $d.Content
$find2.Find.ClearFormatting()
$find2.Find.Replacement.ClearFormatting()
$find2.Find.Execute("May be included as separate pdf", $true, $false, $false, $false, $false, $true, 1, $false, "Included as separate PDF", 2) | Out-Null

# ---------------------------------------------------------------------
# 5) Drop the trailing empty paragraph at the very end of the body (the
#    body always keeps a final paragraph mark, so merge the one before
#    it away instead of trying to delete the true last mark).
# ---------------------------------------------------------------------
$storyEnd = $d.Content.End
$trailingMark = $d.Range($storyEnd - 2, $storyEnd - 1)
if ($trailingMark.Text -eq "`r") {
    $trailingMark.Delete()
}
